# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" sheet (with its fund-holding detail table) right
# after the "总计" summary sheet, and adds a corresponding summary row on
# "总计" itself.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet: insert a new first data row for
#    2022-Q4 and push the existing rows down (re-numbering the index
#    column A as we go).
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

$totals.Range("A5").Value = 3
$totals.Range("B5").Value = "2021-Q1"
$totals.Range("C5").Value = 2
$totals.Range("D5").Value = 0.02

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2021-Q2"
$totals.Range("C4").Value = 1
$totals.Range("D4").Value = 0

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 3
$totals.Range("D3").Value = 0.08

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 5
$totals.Range("D2").Value = 0.44

# carry the index-column style (s="2") down onto the newly added row 5
$totals.Range("A2").Copy()
$totals.Range("A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" sheet right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $totals)
$q4.Name = "2022-Q4"

# Header row
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# copy the bold/centered/bordered header style from the totals sheet
$totals.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# Cell values that look numeric (fund codes, sizes, ratios, weights) must
# stay text, matching the source data - pre-format those cells as Text so
# assigning the literal string doesn't get reinterpreted as a number.
$q4.Range("B2:B6").NumberFormat = "@"
$q4.Range("D2:G6").NumberFormat = "@"

$data = @(
  @(0, "005313", "万家中证1000指数增强A", "22.07", "94.13", "1.03", "0.2273", 5),
  @(1, "005314", "万家中证1000指数增强C", "19.61", "94.13", "1.03", "0.2020", 5),
  @(2, "006354", "国泰民裕进取灵活配置混合", "0.50", "69.40", "1.70", "0.0085", 6),
  @(3, "005095", "国泰量化成长优选混合A", "0.20", "87.25", "2.19", "0.0044", 5),
  @(4, "005096", "国泰量化成长优选混合C", "0.02", "87.25", "2.19", "0.0004", 5)
)

$r = 2
foreach ($row in $data) {
    $q4.Range("A$r").Value = $row[0]
    $q4.Range("B$r").Value = $row[1]
    $q4.Range("C$r").Value = $row[2]
    $q4.Range("D$r").Value = $row[3]
    $q4.Range("E$r").Value = $row[4]
    $q4.Range("F$r").Value = $row[5]
    $q4.Range("G$r").Value = $row[6]
    $q4.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# copy the index-column style (s="2") from the totals sheet onto A2:A6
$totals.Range("A2").Copy()
$q4.Range("A2:A6").PasteSpecial(-4122)
